# "Repayment schedule" sheet gains a new (currently blank) column between
# the existing "In Advance" (M) and "Late" (N) columns - i.e. a brand new,
# empty column N is inserted, pushing the old N (Late), O (heading/"Original"),
# P (Outstanding) columns one slot to the right (-> O, P, Q).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Repayment schedule")

# Insert a new blank column at N - everything from N onward (old Late /
# heading / Outstanding columns) shifts right to O / P / Q.
$ws.Columns("N").Insert()

# The newly inserted column should look like its left neighbour (M, "In
# Advance") width-wise.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab (it moves from
# "NewLoanInput"), with a fresh selection.
$ws.Activate()
$ws.Range("S4").Select()
